$wb = $excel.ActiveWorkbook

# Rename/reposition: add a new first sheet named "Sheet1"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet1"
$newSheet.Move($wb.Worksheets.Item(1))
$ws1 = $wb.Worksheets.Item("Sheet1")

# Fill test data
$data = @(
  @("TestCaseName","MethodName","Value"),
  @("TC_AddNewResource","Uname","vendor2_appian"),
  @("TC_AddNewResource","Password","Welcome21"),
  @("TC_AddNewResource","SOW","GS-16075-SOW-2.00"),
  @("TC_AddNewPosition","Uname","vendor2_appian"),
  @("TC_AddNewPosition","Password","Welcome21"),
  @("TC_AddNewPosition","DemandID",51745),
  @("TC_OrderOwner_Approval","Uname","Order_Owner_01"),
  @("TC_OrderOwner_Approval","Password","Welcome22"),
  @("TC_OrderOwner_Approval","DemandID",51745)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws1.Cells.Item($r, 1).Value = $data[$i][0]
    $ws1.Cells.Item($r, 2).Value = $data[$i][1]
    $ws1.Cells.Item($r, 3).Value = $data[$i][2]
}

# number format for the two DemandID numeric cells
$ws1.Range("C7").NumberFormat = "0"
$ws1.Range("C10").NumberFormat = "0"

# text format for the SOW value cell
$ws1.Range("C4").NumberFormat = "@"

# Hyperlink style (without actual hyperlink) for Uname/Uname-password rows
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Range("C5").Style = "Hyperlink"
$ws1.Range("C8").Style = "Hyperlink"

# Real hyperlinks for the two "Welcome21" password cells
$ws1.Hyperlinks.Add($ws1.Range("C3"), "mailto:Aug@20182018", "", "", "Aug@20182018")
$ws1.Range("C3").Value = "Welcome21"
$ws1.Range("C3").Style = "Hyperlink"

$ws1.Hyperlinks.Add($ws1.Range("C6"), "mailto:Aug@20182018", "", "", "Aug@20182018")
$ws1.Range("C6").Value = "Welcome21"
$ws1.Range("C6").Style = "Hyperlink"

# column widths (autofit based on content)
$ws1.Range("A1:C10").EntireColumn.AutoFit()

# selection + active tab
$ws1.Range("A8").Select()
$ws1.Activate()

# --- admin sheet (now 2nd) ---
$wsAdmin = $wb.Worksheets.Item("admin")
$wsAdmin.Range("B3").Select()

# --- Demand details sheet (now 3rd) ---
$wsDemand = $wb.Worksheets.Item("Demand details")
$wsDemand.Range("A2").Select()

# re-activate Sheet1 so it's the saved active tab
$ws1.Activate()
$ws1.Range("A8").Select()
